$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "https://www.rfebm.com/competiciones/resultados_completos.php?seleccion=0&id=1028265"
$ws.Hyperlinks.Add($ws.Range("A11"), "https://www.rfebm.com/competiciones/resultados_completos.php?seleccion=0&id=1028265")

$ws.Range("A11").Select() | Out-Null
